$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the G9 cell (Derek row) value 1 -> 0
$ws.Range("G9").Value = 0

# Add the "Bench2" shared string first so it is allocated before "Bench1"
# in the shared string table, matching the target workbook ordering.
$ws.Range("A13").Value = "Bench2"

# Update row 12: was "BenchWarmer" with all 4's, now "Bench1" with new values
$ws.Range("A12").Value = "Bench1"
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 3
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0

# Finish populating new row 13: "Bench2"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 2
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0

# Update selection to match the final active cell in the diff
$ws.Range("B12").Select()
